$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with freshly scraped values.
# Values are written as text (matching the source sheet's inlineStr cells),
# using a Text number format while writing so Excel does not reinterpret
# numeric-looking / percent-looking strings as real numbers, then the style
# is reset back to Normal so no formatting residue is left behind.
$updates = @{
    'D2' = '288.32'
    'E2' = '-1.25%'
    'E3' = '1.20%'
    'D4' = '4.927'
    'E4' = '-0.68%'
    'D5' = '0.07346'
    'E5' = '1.61%'
    'D6' = '2.199'
    'E6' = '22.76%'
    'D7' = '7.726'
    'E7' = '0.55%'
    'D8' = '3.734'
    'E8' = '-0.74%'
    'D9' = '0.9027'
    'E9' = '0.54%'
    'D10' = '0.09185'
    'E10' = '19.53%'
    'D11' = '0.1681'
    'E11' = '1.10%'
    'D12' = '0.08243'
    'E12' = '1.87%'
    'D13' = '0.03120'
    'E13' = '2.63%'
    'D14' = '0.09946'
    'E14' = '-0.73%'
    'D15' = '0.001498'
    'E15' = '-0.19%'
    'D16' = '0.005701'
    'E16' = '-0.89%'
    'E17' = '1.46%'
    'D18' = '2.058'
    'E18' = '-1.25%'
    'D19' = '0.3332'
    'E19' = '0.41%'
    'D20' = '0.1293'
    'E20' = '0.29%'
    'D21' = '4.154'
    'E21' = '2.76%'
    'D22' = '0.2104'
    'E22' = '-6.62%'
    'D23' = '0.04537'
    'E23' = '0.66%'
    'E24' = '-0.63%'
    'D25' = '0.004152'
    'E25' = '3.46%'
    'E26' = '4.09%'
    'D27' = '0.0003399'
    'D39' = '0.01561'
    'E39' = '-2.29%'
    'D40' = '0.04437'
    'E40' = '0.27%'
    'D41' = '0.007407'
    'E41' = '2.23%'
    'D42' = '0.008991'
    'E42' = '17.05%'
    'D43' = '0.1324'
    'E43' = '1.10%'
    'D44' = '0.002231'
    'E44' = '10.46%'
    'D45' = '0.009111'
    'E45' = '-4.11%'
    'D46' = '0.00006123'
    'E46' = '3.28%'
    'E47' = '0.16%'
    'D48' = '2.318'
    'E48' = '3.24%'
    'D49' = '0.002003'
    'E49' = '-33.33%'
    'D50' = '0.00002103'
    'E50' = '0.16%'
    'D51' = '0.0002003'
    'E51' = '0.16%'
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$cellRef]
    $cell.Style = "Normal"
}
